$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.999.98'
$ws.Range('E2').Value = '  -0.52%  '

$ws.Range('D3').Value = '1.639.69'
$ws.Range('E3').Value = '  -1.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.34%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5052'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.16%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06452'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.14%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2574'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E10').Value = '  -1.92%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07700'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.00%  '

$ws.Range('D12').Value = '1.641.27'
$ws.Range('E12').Value = '  -0.90%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.247'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.17%  '

$ws.Range('D14').Value = '1.865.80'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5451'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.42%  '

$ws.Range('D16').Value = '0.0₅7934'
$ws.Range('E16').Value = '  -1.40%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.47%  '

$ws.Range('D18').Value = '25.987.47'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.288'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.23%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.994'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.973'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.961'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.43%  '

$ws.Range('E27').Value = '  -0.73%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.10%  '

$ws.Range('E29').Value = '  -3.58%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05049'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.60%  '

$ws.Range('E31').Value = '  -0.95%  '

$ws.Range('E32').Value = '  -3.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.194'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.21%  '

$ws.Range('E34').Value = '  -1.94%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.343'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.640'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.38%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8925'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.24%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5623'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.17%  '

$ws.Range('D39').Value = '1.148.27'
$ws.Range('E39').Value = '  -1.22%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01570'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.45%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.560'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.36%  '

$ws.Range('E42').Value = '  +0.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.664'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.27%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8079'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.51%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.38%  '

$ws.Range('D46').Value = '1.777.60'

$ws.Range('D47').Value = '0.0₈113'
$ws.Range('E47').Value = '  -0.56%  '

$ws.Range('E48').Value = '  +0.22%  '

$ws.Range('E49').Value = '  +0.00%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '55.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.65%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05033'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.80%  '

